$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 55 (J05AJ03 / Dolutegravir entry),
# shifting all subsequent rows down by 2.
$ws.Rows("55:56").Insert()

# Row 55: J05AJ01 / Raltegravir / RAL
$ws.Cells.Item(55, 1).Value = "J05AJ01 "
$ws.Cells.Item(55, 2).Value = "Raltegravir "
$ws.Cells.Item(55, 3).Value = "II"
$ws.Cells.Item(55, 4).Value = "RAL"
$ws.Cells.Item(55, 5).Value = "RAL"
$ws.Cells.Item(55, 6).Value = "II"

# Row 56: J05AJ02 / Elvitegravir / EVG
$ws.Cells.Item(56, 1).Value = "J05AJ02"
$ws.Cells.Item(56, 2).Value = "Elvitegravir"
$ws.Cells.Item(56, 3).Value = "II"
$ws.Cells.Item(56, 4).Value = "EVG"
$ws.Cells.Item(56, 5).Value = "EVG"
$ws.Cells.Item(56, 6).Value = "II"

# Colour the two new rows red (same style used for other newly-added entries).
$ws.Range("A55:F56").Font.Color = 255

# Match the saved selection state.
$ws.Range("A30").Select() | Out-Null
